$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B5 (date) and E5 (CDC) values; C5/F5 recalc via formulas.
$ws.Range("B5").Value = 44533
$ws.Range("E5").Value = 7534

# F5 becomes a standalone (non-shared) formula matching the others' pattern.
$ws.Range("F5").Formula = '=IF(E5 > 0, DATE(2001,5,1)+E5-1, "cdc inválido")'

# Update the active selection on the sheet.
$ws.Range("J9").Select()
